$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Status column (C) for rows 16 and 23 from "In Progress" to "Done"
$ws.Range("C16").Value = "Done"
$ws.Range("C23").Value = "Done"
